$d = $word.ActiveDocument

# Add the two new custom (quick-style) paragraph styles, each based on
# "Normal", matching the "acpt: add scenarios for style access" commit:
#   - "Foo Bar" (styleId FooBar)
#   - "Bar Foo" (styleId BarFoo)

$normal = $d.Styles.Item("Normal")

$fooBar = $d.Styles.Add("Foo Bar", 1)
$fooBar.BaseStyle = $normal
$fooBar.QuickStyle = $true

$barFoo = $d.Styles.Add("Bar Foo", 1)
$barFoo.BaseStyle = $normal
$barFoo.QuickStyle = $true
